$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows before row 124 (shifts old rows 124-132 down to 127-135) ---
$ws.Rows.Item(124).Insert()
$ws.Rows.Item(124).Insert()
$ws.Rows.Item(124).Insert()

# Common (constant) values shared by every Chirimoya / Comercializadora del Agro de Limari row
$commonA = 2
$commonB = "Comercializadora del Agro de Limarí"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100107
$commonH = "Otros"
$commonI = 100107002
$commonJ = "Chirimoya"
$commonK = "Cultivar IV Región"
$commonR = "Provincia de Limarí"

# --- Rows 121-123: updated with the new (2022-11-24) weekly report values ---
$ws.Cells.Item(121,1).Value = $commonA
$ws.Cells.Item(121,2).Value = $commonB
$ws.Cells.Item(121,3).Value = $commonC
$ws.Cells.Item(121,4).Value = 44889
$ws.Cells.Item(121,5).Value = $commonE
$ws.Cells.Item(121,6).Value = $commonF
$ws.Cells.Item(121,7).Value = $commonG
$ws.Cells.Item(121,8).Value = $commonH
$ws.Cells.Item(121,9).Value = $commonI
$ws.Cells.Item(121,10).Value = $commonJ
$ws.Cells.Item(121,11).Value = $commonK
$ws.Cells.Item(121,12).Value = "Especial"
$ws.Cells.Item(121,13).Value = 300
$ws.Cells.Item(121,14).Value = 18000
$ws.Cells.Item(121,15).Value = 19000
$ws.Cells.Item(121,16).Value = 18500
$ws.Cells.Item(121,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(121,18).Value = $commonR
$ws.Cells.Item(121,19).Value = 1850
$ws.Cells.Item(121,20).Value = 10

$ws.Cells.Item(122,1).Value = $commonA
$ws.Cells.Item(122,2).Value = $commonB
$ws.Cells.Item(122,3).Value = $commonC
$ws.Cells.Item(122,4).Value = 44889
$ws.Cells.Item(122,5).Value = $commonE
$ws.Cells.Item(122,6).Value = $commonF
$ws.Cells.Item(122,7).Value = $commonG
$ws.Cells.Item(122,8).Value = $commonH
$ws.Cells.Item(122,9).Value = $commonI
$ws.Cells.Item(122,10).Value = $commonJ
$ws.Cells.Item(122,11).Value = $commonK
$ws.Cells.Item(122,12).Value = "Primera"
$ws.Cells.Item(122,13).Value = 240
$ws.Cells.Item(122,14).Value = 15000
$ws.Cells.Item(122,15).Value = 16000
$ws.Cells.Item(122,16).Value = 15500
$ws.Cells.Item(122,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(122,18).Value = $commonR
$ws.Cells.Item(122,19).Value = 1550
$ws.Cells.Item(122,20).Value = 10

$ws.Cells.Item(123,1).Value = $commonA
$ws.Cells.Item(123,2).Value = $commonB
$ws.Cells.Item(123,3).Value = $commonC
$ws.Cells.Item(123,4).Value = 44889
$ws.Cells.Item(123,5).Value = $commonE
$ws.Cells.Item(123,6).Value = $commonF
$ws.Cells.Item(123,7).Value = $commonG
$ws.Cells.Item(123,8).Value = $commonH
$ws.Cells.Item(123,9).Value = $commonI
$ws.Cells.Item(123,10).Value = $commonJ
$ws.Cells.Item(123,11).Value = $commonK
$ws.Cells.Item(123,12).Value = "Segunda"
$ws.Cells.Item(123,13).Value = 200
$ws.Cells.Item(123,14).Value = 11000
$ws.Cells.Item(123,15).Value = 12000
$ws.Cells.Item(123,16).Value = 11500
$ws.Cells.Item(123,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(123,18).Value = $commonR
$ws.Cells.Item(123,19).Value = 1150
$ws.Cells.Item(123,20).Value = 10

# --- New rows 124-126: the previous (2020-11-26) weekly report values that used to live in 121-123 ---
$ws.Cells.Item(124,1).Value = $commonA
$ws.Cells.Item(124,2).Value = $commonB
$ws.Cells.Item(124,3).Value = $commonC
$ws.Cells.Item(124,4).Value = 44161
$ws.Cells.Item(124,5).Value = $commonE
$ws.Cells.Item(124,6).Value = $commonF
$ws.Cells.Item(124,7).Value = $commonG
$ws.Cells.Item(124,8).Value = $commonH
$ws.Cells.Item(124,9).Value = $commonI
$ws.Cells.Item(124,10).Value = $commonJ
$ws.Cells.Item(124,11).Value = $commonK
$ws.Cells.Item(124,12).Value = "Especial"
$ws.Cells.Item(124,13).Value = 240
$ws.Cells.Item(124,14).Value = 13000
$ws.Cells.Item(124,15).Value = 13500
$ws.Cells.Item(124,16).Value = 13250
$ws.Cells.Item(124,17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(124,18).Value = $commonR
$ws.Cells.Item(124,19).Value = 1656
$ws.Cells.Item(124,20).Value = 8

$ws.Cells.Item(125,1).Value = $commonA
$ws.Cells.Item(125,2).Value = $commonB
$ws.Cells.Item(125,3).Value = $commonC
$ws.Cells.Item(125,4).Value = 44161
$ws.Cells.Item(125,5).Value = $commonE
$ws.Cells.Item(125,6).Value = $commonF
$ws.Cells.Item(125,7).Value = $commonG
$ws.Cells.Item(125,8).Value = $commonH
$ws.Cells.Item(125,9).Value = $commonI
$ws.Cells.Item(125,10).Value = $commonJ
$ws.Cells.Item(125,11).Value = $commonK
$ws.Cells.Item(125,12).Value = "Primera"
$ws.Cells.Item(125,13).Value = 240
$ws.Cells.Item(125,14).Value = 11000
$ws.Cells.Item(125,15).Value = 11500
$ws.Cells.Item(125,16).Value = 11250
$ws.Cells.Item(125,17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(125,18).Value = $commonR
$ws.Cells.Item(125,19).Value = 1406
$ws.Cells.Item(125,20).Value = 8

$ws.Cells.Item(126,1).Value = $commonA
$ws.Cells.Item(126,2).Value = $commonB
$ws.Cells.Item(126,3).Value = $commonC
$ws.Cells.Item(126,4).Value = 44161
$ws.Cells.Item(126,5).Value = $commonE
$ws.Cells.Item(126,6).Value = $commonF
$ws.Cells.Item(126,7).Value = $commonG
$ws.Cells.Item(126,8).Value = $commonH
$ws.Cells.Item(126,9).Value = $commonI
$ws.Cells.Item(126,10).Value = $commonJ
$ws.Cells.Item(126,11).Value = $commonK
$ws.Cells.Item(126,12).Value = "Segunda"
$ws.Cells.Item(126,13).Value = 200
$ws.Cells.Item(126,14).Value = 9000
$ws.Cells.Item(126,15).Value = 9500
$ws.Cells.Item(126,16).Value = 9250
$ws.Cells.Item(126,17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(126,18).Value = $commonR
$ws.Cells.Item(126,19).Value = 1156
$ws.Cells.Item(126,20).Value = 8

# Apply the date number format (style used elsewhere in column D) to the D cells of the newly
# inserted rows, matching the style already present on D121:D123 before the insert.
$ws.Range("D124:D126").NumberFormat = $ws.Range("D127").NumberFormat
